$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 183, shifting existing rows 183:197 down to 184:198
$ws.Rows.Item(183).Insert()

# Populate the newly inserted row 183 with the new weekly record
$ws.Range("A183").Value = 3
$ws.Range("B183").Value = "Femacal de La Calera"
$ws.Range("C183").Value = "Coquimbo"
$ws.Range("D183").Value = 44461
$ws.Range("E183").Value = 5
$ws.Range("F183").Value = 100112043
$ws.Range("G183").Value = "Pepino ensalada"
$ws.Range("H183").Value = "Sin especificar"
$ws.Range("I183").Value = "Primera"
$ws.Range("J183").Value = 80
$ws.Range("K183").Value = 14000
$ws.Range("L183").Value = 15000
$ws.Range("M183").Value = 14500
$ws.Range("N183").Value = "`$/caja 70 unidades"
$ws.Range("O183").Value = "Región de Arica y Parinacota"
$ws.Range("P183").Value = 207
$ws.Range("Q183").Value = 70
$ws.Range("R183").Value = "Hortaliza"
